$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd21fab4be0>),
                (''model'',
                 BaggingClassifier(estimator=RandomForestClassifier(max_depth=1,
                                                                    min_samples_leaf=4,
                                                                    min_samples_split=9,
                                                                    n_estimators=50,
                                                                    random_state=42),
                                   random_state=42))])'
$ws.Range("B2").Value = 0.7452564102564103
$ws.Range("C2").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd15c47c220>, ''scaler'': None, ''model__n_estimators'': 10, ''model__estimator__n_estimators'': 50, ''model__estimator__min_samples_split'': 9, ''model__estimator__min_samples_leaf'': 4, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 1, ''model__estimator__class_weight'': None}'
$ws.Range("D2").Value = 0.7792411300436752
$ws.Range("E2").Value = 0.6426895659895661
$ws.Range("F2").Value = 0.8108108108108109
$ws.Range("G2").Value = 0.745020423207164
$ws.Range("H2").Value = 0.6079535714285714
$ws.Range("I2").Value = 0.7142857142857143
$ws.Range("J2").Value = 0.8370000000000002
$ws.Range("K2").Value = 0.7181666666666667
$ws.Range("N2").Value = '[1 1 1 1 1 1 0 1 1 1 1 1 1 1 0 0 1 1 1 1 1 1 1 1]'

$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd15c7112b0>),
                (''model'',
                 BaggingClassifier(estimator=RandomForestClassifier(max_depth=1,
                                                                    min_samples_leaf=11,
                                                                    min_samples_split=7,
                                                                    n_estimators=10,
                                                                    random_state=42),
                                   random_state=42))])'
$ws.Range("B3").Value = 0.7402930402930403
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd15c4fbd90>, ''scaler'': StandardScaler(), ''model__n_estimators'': 10, ''model__estimator__n_estimators'': 10, ''model__estimator__min_samples_split'': 7, ''model__estimator__min_samples_leaf'': 11, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 1, ''model__estimator__class_weight'': None}'
$ws.Range("D3").Value = 0.7800628670967936
$ws.Range("E3").Value = 0.616917882117882
$ws.Range("G3").Value = 0.7497937816197513
$ws.Range("H3").Value = 0.6661603174603176
$ws.Range("J3").Value = 0.8310425531914896
$ws.Range("K3").Value = 0.6243333333333334

$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 BaggingClassifier(estimator=RandomForestClassifier(max_depth=1,
                                                                    max_features=''log2'',
                                                                    min_samples_leaf=10,
                                                                    n_estimators=50,
                                                                    random_state=42),
                                   random_state=42))])'
$ws.Range("B4").Value = 0.7256010656010656
$ws.Range("C4").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': StandardScaler(), ''model__n_estimators'': 10, ''model__estimator__n_estimators'': 50, ''model__estimator__min_samples_split'': 2, ''model__estimator__min_samples_leaf'': 10, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 1, ''model__estimator__class_weight'': None}'
$ws.Range("D4").Value = 0.759687568549166
$ws.Range("E4").Value = 0.6224345876345875
$ws.Range("F4").Value = 0.7222222222222222
$ws.Range("G4").Value = 0.7427976912452163
$ws.Range("H4").Value = 0.6294428571428571
$ws.Range("I4").Value = 0.7647058823529411
$ws.Range("J4").Value = 0.7896000000000001
$ws.Range("K4").Value = 0.6512
$ws.Range("L4").Value = 0.6842105263157895
$ws.Range("N4").Value = '[1 1 1 1 1 1 0 1 1 0 1 0 0 1 1 1 0 1 1 0 1 0 1 1]'

$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd15c54a280>),
                (''model'',
                 BaggingClassifier(estimator=RandomForestClassifier(max_depth=1,
                                                                    max_features=''log2'',
                                                                    min_samples_leaf=4,
                                                                    min_samples_split=4,
                                                                    random_state=42),
                                   random_state=42))])'
$ws.Range("B5").Value = 0.7578571428571428
$ws.Range("C5").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd15c0d9400>, ''scaler'': StandardScaler(), ''model__n_estimators'': 10, ''model__estimator__n_estimators'': 100, ''model__estimator__min_samples_split'': 4, ''model__estimator__min_samples_leaf'': 4, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 1, ''model__estimator__class_weight'': None}'
$ws.Range("D5").Value = 0.7759532565515141
$ws.Range("E5").Value = 0.6526460539460538
$ws.Range("G5").Value = 0.713691858899535
$ws.Range("H5").Value = 0.6084265873015873
$ws.Range("J5").Value = 0.8740816326530613
$ws.Range("K5").Value = 0.7476666666666668
$ws.Range("N5").Value = '[1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'

$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd15c54a100>),
                (''model'',
                 BaggingClassifier(estimator=RandomForestClassifier(max_depth=1,
                                                                    max_features=''log2'',
                                                                    min_samples_leaf=9,
                                                                    min_samples_split=3,
                                                                    n_estimators=10,
                                                                    random_state=42),
                                   n_estimators=50, random_state=42))])'
$ws.Range("B6").Value = 0.7400183150183149
$ws.Range("C6").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd15c50a6a0>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 50, ''model__estimator__n_estimators'': 10, ''model__estimator__min_samples_split'': 3, ''model__estimator__min_samples_leaf'': 9, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 1, ''model__estimator__class_weight'': None}'
$ws.Range("D6").Value = 0.7981627198965039
$ws.Range("E6").Value = 0.6560632423132422
$ws.Range("F6").Value = 0.6470588235294118
$ws.Range("G6").Value = 0.7747333667104463
$ws.Range("H6").Value = 0.6223992063492063
$ws.Range("I6").Value = 0.4782608695652174
$ws.Range("J6").Value = 0.845153846153846
$ws.Range("K6").Value = 0.7278333333333334
$ws.Range("L6").Value = 1
$ws.Range("N6").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1]'

Write-Host "Applied edits"